$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New log entry row (row 36)
$row = 36

$ws.Cells.Item($row, 1).Value2 = "Verzoek om factuur"
$ws.Cells.Item($row, 2).Value2 = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value2 = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$ws.Cells.Item($row, 4).Value2 = "Factuur / Administratie"
$ws.Cells.Item($row, 5).Value2 = "Beste klant,`r`n`r`nBedankt voor uw e-mail. Om u te kunnen helpen met het versturen van een factuur voor uw laatste bestelling, hebben wij wat extra informatie nodig. Kunt u ons alstublieft de volgende gegevens verstrekken:`r`n`r`n1. Uw bestelnummer`r`n2. De datum van uw bestelling`r`n3. Het e-mailadres waarnaar wij de factuur kunnen sturen`r`n`r`nZodra wij deze gegevens van u hebben ontvangen, zullen wij zo spoedig mogelijk de factuur voor u opstellen en toesturen.`r`n`r`nMet vriendelijke groet,`r`n`r`n[Bedrijfsnaam] E-mailassistent"
$ws.Cells.Item($row, 6).Value2 = "2025-06-24 21:43:40"
$ws.Cells.Item($row, 7).Value2 = "Ja"

# Prevent the engine's auto row-height recalculation (triggered by the
# embedded line breaks) from stamping an explicit customHeight on the row.
$ws.Rows.Item($row).AutoFit()

# Extend the existing conditional formatting rules so they keep covering
# the Categorie (D) and Beantwoord (G) columns through the new row.
$fcsD = $ws.Range("D2:D35").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D36"))
}

$fcsG = $ws.Range("G2:G35").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G36"))
}

# Update the Dashboard summary count for "Factuur / Administratie".
$wsd = $wb.Worksheets.Item("Dashboard")
$wsd.Cells.Item(3, 2).Value2 = 5
